# Remove the "ocid" column (column A) from every sheet except "Activity".
# This shifts all remaining columns one position to the left on each of
# those sheets and drops the now-unused "ocid" shared string on save.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Classification", "Documents", "Event", "GrantProgramme", "Location", "Organization", "Transaction")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(1).Delete()
}
